$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "C3"
$ws.Range("C2").Value2 = "Itgax"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1.468507333333333
$ws.Range("H2").Value2 = 4.405521999999999
$ws.Range("I2").Value2 = 0.005118279455112885
$ws.Range("J2").Value2 = 0.005118279455112885
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.05829066666666666
$ws.Range("N2").Value2 = 0.174872
$ws.Range("O2").Value2 = 0.0009372933078703915
$ws.Range("P2").Value2 = 0.0009372933078703916
$ws.Range("Q2").Value2 = 0.08560027146488887
$ws.Range("R2").Value2 = 0.7704024431839999
$ws.Range("S2").Value2 = 0.000004797329081087821
$ws.Range("T2").Value2 = 0.000004797329081087821

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "C3"
$ws.Range("C3").Value2 = "Itgax"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1.468507333333333
$ws.Range("H3").Value2 = 4.405521999999999
$ws.Range("I3").Value2 = 0.005118279455112885
$ws.Range("J3").Value2 = 0.005118279455112885
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.01197066666666667
$ws.Range("N3").Value2 = 0.035912
$ws.Range("O3").Value2 = 0.0001924840870593434
$ws.Range("P3").Value2 = 0.0001924840870593434
$ws.Range("Q3").Value2 = 0.01757901178488889
$ws.Range("R3").Value2 = 0.158211106064
$ws.Range("S3").Value2 = 0.0000009851873482319972
$ws.Range("T3").Value2 = 0.0000009851873482319972

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "C3"
$ws.Range("C4").Value2 = "Itgax"
$ws.Range("D4").Value2 = "Resolving-Mac"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1.468507333333333
$ws.Range("H4").Value2 = 4.405521999999999
$ws.Range("I4").Value2 = 0.005118279455112885
$ws.Range("J4").Value2 = 0.005118279455112885
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 62.120161
$ws.Range("N4").Value2 = 186.360483
$ws.Range("O4").Value2 = 0.9988702226050702
$ws.Range("P4").Value2 = 0.9988702226050703
$ws.Range("Q4").Value2 = 91.22391197634731
$ws.Range("R4").Value2 = 821.0152077871259
$ws.Range("S4").Value2 = 0.005112496938683564
$ws.Range("T4").Value2 = 0.005112496938683565

# Row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "C3"
$ws.Range("C5").Value2 = "Itgax"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 259.5505726666667
$ws.Range("H5").Value2 = 778.6517180000001
$ws.Range("I5").Value2 = 0.9046276674881553
$ws.Range("J5").Value2 = 0.9046276674881553
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.05829066666666666
$ws.Range("N5").Value2 = 0.174872
$ws.Range("O5").Value2 = 0.0009372933078703915
$ws.Range("P5").Value2 = 0.0009372933078703916
$ws.Range("Q5").Value2 = 15.12937591445511
$ws.Range("R5").Value2 = 136.164383230096
$ws.Range("S5").Value2 = 0.0008479014588510497
$ws.Range("T5").Value2 = 0.0008479014588510497

# Row 6
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "C3"
$ws.Range("C6").Value2 = "Itgax"
$ws.Range("D6").Value2 = "MuSCs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 259.5505726666667
$ws.Range("H6").Value2 = 778.6517180000001
$ws.Range("I6").Value2 = 0.9046276674881553
$ws.Range("J6").Value2 = 0.9046276674881553
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.01197066666666667
$ws.Range("N6").Value2 = 0.035912
$ws.Range("O6").Value2 = 0.0001924840870593434
$ws.Range("P6").Value2 = 0.0001924840870593434
$ws.Range("Q6").Value2 = 3.106993388535111
$ws.Range("R6").Value2 = 27.962940496816
$ws.Range("S6").Value2 = 0.0001741264307050808
$ws.Range("T6").Value2 = 0.0001741264307050808

# Row 7
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "C3"
$ws.Range("C7").Value2 = "Itgax"
$ws.Range("D7").Value2 = "Resolving-Mac"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 259.5505726666667
$ws.Range("H7").Value2 = 778.6517180000001
$ws.Range("I7").Value2 = 0.9046276674881553
$ws.Range("J7").Value2 = 0.9046276674881553
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 62.120161
$ws.Range("N7").Value2 = 186.360483
$ws.Range("O7").Value2 = 0.9988702226050702
$ws.Range("P7").Value2 = 0.9988702226050703
$ws.Range("Q7").Value2 = 16123.32336169553
$ws.Range("R7").Value2 = 145109.9102552598
$ws.Range("S7").Value2 = 0.9036056395985991
$ws.Range("T7").Value2 = 0.9036056395985992

# Row 8
$ws.Range("A8").Value2 = "MuSCs"
$ws.Range("B8").Value2 = "C3"
$ws.Range("C8").Value2 = "Itgax"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 0.5890733333333333
$ws.Range("H8").Value2 = 1.76722
$ws.Range("I8").Value2 = 0.002053133730501083
$ws.Range("J8").Value2 = 0.002053133730501083
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.05829066666666666
$ws.Range("N8").Value2 = 0.174872
$ws.Range("O8").Value2 = 0.0009372933078703915
$ws.Range("P8").Value2 = 0.0009372933078703916
$ws.Range("Q8").Value2 = 0.03433747731555555
$ws.Range("R8").Value2 = 0.30903729584
$ws.Range("S8").Value2 = 0.000001924388505761637
$ws.Range("T8").Value2 = 0.000001924388505761637

# Row 9
$ws.Range("A9").Value2 = "MuSCs"
$ws.Range("B9").Value2 = "C3"
$ws.Range("C9").Value2 = "Itgax"
$ws.Range("D9").Value2 = "MuSCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 0.5890733333333333
$ws.Range("H9").Value2 = 1.76722
$ws.Range("I9").Value2 = 0.002053133730501083
$ws.Range("J9").Value2 = 0.002053133730501083
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.01197066666666667
$ws.Range("N9").Value2 = 0.035912
$ws.Range("O9").Value2 = 0.0001924840870593434
$ws.Range("P9").Value2 = 0.0001924840870593434
$ws.Range("Q9").Value2 = 0.007051600515555556
$ws.Range("R9").Value2 = 0.06346440464
$ws.Range("S9").Value2 = 0.000000395195571726245
$ws.Range("T9").Value2 = 0.000000395195571726245

# Row 10
$ws.Range("A10").Value2 = "MuSCs"
$ws.Range("B10").Value2 = "C3"
$ws.Range("C10").Value2 = "Itgax"
$ws.Range("D10").Value2 = "Resolving-Mac"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 0.5890733333333333
$ws.Range("H10").Value2 = 1.76722
$ws.Range("I10").Value2 = 0.002053133730501083
$ws.Range("J10").Value2 = 0.002053133730501083
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 62.120161
$ws.Range("N10").Value2 = 186.360483
$ws.Range("O10").Value2 = 0.9988702226050702
$ws.Range("P10").Value2 = 0.9988702226050703
$ws.Range("Q10").Value2 = 36.59333030747333
$ws.Range("R10").Value2 = 329.33997276726
$ws.Range("S10").Value2 = 0.002050814146423595
$ws.Range("T10").Value2 = 0.002050814146423595

# Row 11
$ws.Range("A11").Value2 = "Resolving-Mac"
$ws.Range("B11").Value2 = "C3"
$ws.Range("C11").Value2 = "Itgax"
$ws.Range("D11").Value2 = "ECs"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 25.306101
$ws.Range("H11").Value2 = 75.91830299999999
$ws.Range("I11").Value2 = 0.0882009193262308
$ws.Range("J11").Value2 = 0.0882009193262308
$ws.Range("K11").Value2 = 1
$ws.Range("L11").Value2 = 0.3333333333333333
$ws.Range("M11").Value2 = 0.05829066666666666
$ws.Range("N11").Value2 = 0.174872
$ws.Range("O11").Value2 = 0.0009372933078703915
$ws.Range("P11").Value2 = 0.0009372933078703916
$ws.Range("Q11").Value2 = 1.475109498024
$ws.Range("R11").Value2 = 13.275985482216
$ws.Range("S11").Value2 = 0.00008267013143249241
$ws.Range("T11").Value2 = 0.00008267013143249241

# Row 12
$ws.Range("A12").Value2 = "Resolving-Mac"
$ws.Range("B12").Value2 = "C3"
$ws.Range("C12").Value2 = "Itgax"
$ws.Range("D12").Value2 = "MuSCs"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 25.306101
$ws.Range("H12").Value2 = 75.91830299999999
$ws.Range("I12").Value2 = 0.0882009193262308
$ws.Range("J12").Value2 = 0.0882009193262308
$ws.Range("K12").Value2 = 1
$ws.Range("L12").Value2 = 0.3333333333333333
$ws.Range("M12").Value2 = 0.01197066666666667
$ws.Range("N12").Value2 = 0.035912
$ws.Range("O12").Value2 = 0.0001924840870593434
$ws.Range("P12").Value2 = 0.0001924840870593434
$ws.Range("Q12").Value2 = 0.3029308997039999
$ws.Range("R12").Value2 = 2.726378097336
$ws.Range("S12").Value2 = 0.00001697727343430433
$ws.Range("T12").Value2 = 0.00001697727343430433

# Row 13
$ws.Range("A13").Value2 = "Resolving-Mac"
$ws.Range("B13").Value2 = "C3"
$ws.Range("C13").Value2 = "Itgax"
$ws.Range("D13").Value2 = "Resolving-Mac"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 25.306101
$ws.Range("H13").Value2 = 75.91830299999999
$ws.Range("I13").Value2 = 0.0882009193262308
$ws.Range("J13").Value2 = 0.0882009193262308
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 62.120161
$ws.Range("N13").Value2 = 186.360483
$ws.Range("O13").Value2 = 0.9988702226050702
$ws.Range("P13").Value2 = 0.9988702226050703
$ws.Range("Q13").Value2 = 1572.019068402261
$ws.Range("R13").Value2 = 14148.17161562035
$ws.Range("S13").Value2 = 0.08810127192136399
$ws.Range("T13").Value2 = 0.08810127192136401

